# "Fix dimes siones file" -- tidy-up pass over dimensiones.xlsx.
#
# 1) The first sheet's tab name had a stray trailing space
#    ("Departamento ") -- trim it to "Departamento".
# 2) Every sheet's header row (row 1) was carrying a stale explicit
#    row height (ht="15.75" / ht="31.5") left over from an older
#    Excel build/font metrics. Auto-fitting the row drops the
#    leftover explicit height so the row falls back to the sheet's
#    default row height again.

$wb = $excel.ActiveWorkbook

# --- 1) Fix the trailing-space typo in the first sheet's name ---
$firstSheet = $wb.Worksheets.Item(1)
if ($firstSheet.Name -ne $firstSheet.Name.Trim()) {
    $firstSheet.Name = $firstSheet.Name.Trim()
}

# --- 2) Re-autofit the header row on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).AutoFit()
}
